# Deploying to gh-pages from @ eurovibes/huibike@27b0568947115e0d296c15c50f3c558930dfdfed
# Updates the BoM sheet ("BoM") with the full per-component breakdown (8 component
# groups instead of the collapsed single "1 uF" row), recolors/resizes accordingly,
# and repositions the KiBot logo picture to match the new sheet extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# ---------------------------------------------------------------------------
# 1. Column widths
# ---------------------------------------------------------------------------
$ws.Columns(3).ColumnWidth = 18.7109375   # C
$ws.Columns(5).ColumnWidth = 30.7109375   # E
$ws.Columns(6).ColumnWidth = 21.7109375   # F
$ws.Columns(7).ColumnWidth = 36.7109375   # G
$ws.Columns(8).ColumnWidth = 60.7109375   # H

# ---------------------------------------------------------------------------
# 2. Component-groups summary value (F2: 1 -> 8)
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = 8

# ---------------------------------------------------------------------------
# 3. Row heights for the new taller BOM rows
# ---------------------------------------------------------------------------
$ws.Range("A9:I9").RowHeight = 30
$ws.Range("A10:I10").RowHeight = 30
$ws.Range("A12:I12").RowHeight = 30
$ws.Range("A14:I14").RowHeight = 30
$ws.Range("A15:I15").RowHeight = 30
$ws.Range("A16:I16").RowHeight = 30

# ---------------------------------------------------------------------------
# 4. Build the new "light" fill styles (used by the even BOM rows 10/12/14/16)
#    Each color is seeded once on its own cell and then fanned out via
#    Copy + PasteSpecial(formats) so the style table stays deduplicated.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

$ws.Range("A10").WrapText = $true
$ws.Range("A10").HorizontalAlignment = 7
$ws.Range("A10").VerticalAlignment = -4108
$ws.Range("A10").Interior.Color = 16056304   # F0FFF4

$ws.Range("B10").WrapText = $true
$ws.Range("B10").HorizontalAlignment = 7
$ws.Range("B10").VerticalAlignment = -4108
$ws.Range("B10").Interior.Color = 12447999   # FFF0BD

$ws.Range("D10").WrapText = $true
$ws.Range("D10").HorizontalAlignment = 7
$ws.Range("D10").VerticalAlignment = -4108
$ws.Range("D10").Interior.Color = 16777200   # F0FFFF

$ws.Range("A10").Copy()
$ws.Range("I10").PasteSpecial($xlPasteFormats)

$ws.Range("B10").Copy()
$ws.Range("C10").PasteSpecial($xlPasteFormats)
$ws.Range("G10").PasteSpecial($xlPasteFormats)

$ws.Range("D10").Copy()
$ws.Range("E10").PasteSpecial($xlPasteFormats)
$ws.Range("F10").PasteSpecial($xlPasteFormats)
$ws.Range("H10").PasteSpecial($xlPasteFormats)

# Fan the now-complete row10 formatting out to the other even rows
$ws.Range("A10:I10").Copy()
$ws.Range("A12:I12").PasteSpecial($xlPasteFormats)
$ws.Range("A14:I14").PasteSpecial($xlPasteFormats)
$ws.Range("A16:I16").PasteSpecial($xlPasteFormats)

# Fan row9's existing (pre-edit) formatting out to the other odd rows
$ws.Range("A9:I9").Copy()
$ws.Range("A11:I11").PasteSpecial($xlPasteFormats)
$ws.Range("A13:I13").PasteSpecial($xlPasteFormats)
$ws.Range("A15:I15").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5. Cell values for the expanded BOM (rows 9-16, one pair of rows per
#    component group: field-name row then field-value row)
# ---------------------------------------------------------------------------

# Row 9 (group 1, field names) / Row 10 (group 1 values) - 100 nF capacitor
$ws.Range("A9").Value = "1"
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "100 nF"
$ws.Range("D9").Value = "CL05B104KO5NNNC"
$ws.Range("E9").Value = "Samsung Electro-Mechanics"
$ws.Range("F9").Value = "C1525"
$ws.Range("G9").Value = "C_0402_1005Metric"
$ws.Range("H9").Value = "16V 100nF X7R ±10% 0402  Multilayer Ceramic Capacitors MLCC - SMD/SMT ROHS"
$ws.Range("I9").Value = "2"

$ws.Range("A10").Value = "2"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "1 µF"
$ws.Range("D10").Value = "CL21B105KBFNNNE"
$ws.Range("E10").Value = "Samsung Electro-Mechanics"
$ws.Range("F10").Value = "C28323"
$ws.Range("G10").Value = "C_0805_2012Metric"
$ws.Range("H10").Value = "50V 1uF X7R ±10% 0805 Multilayer Ceramic Capacitors MLCC - SMD/SMT ROHS"
$ws.Range("I10").Value = "2"

$ws.Range("A11").Value = "3"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "ZMM3V3-M"
$ws.Range("D11").Value = "ZMM3V3-M"
$ws.Range("E11").Value = "ST(Semtech)"
$ws.Range("F11").Value = "C8056"
$ws.Range("G11").Value = "D_MiniMELF"
$ws.Range("H11").Value = "Single 3.1V~3.5V 500mW 3.3V LL-34 Zener Diodes ROHS"
$ws.Range("I11").Value = "2"

$ws.Range("A12").Value = "4"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "Sensor"
$ws.Range("D12").Value = "N/A"
$ws.Range("E12").Value = "N/A"
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = "PinHeader_1x06_P2.54mm_Vertical"
$ws.Range("H12").Value = "Generic connector, single row, 01x06, script generated (kicad-library-utils/schlib/autogen/connector/)"
$ws.Range("I12").Value = "2"

$ws.Range("A13").Value = "5"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "BC847"
$ws.Range("D13").Value = "BC847"
$ws.Range("E13").Value = "Shikues"
$ws.Range("F13").Value = "C475630"
$ws.Range("G13").Value = "SOT-23"
$ws.Range("H13").Value = "45V 200mW 100mA NPN SOT-23 Bipolar Transistors - BJT ROHS"
$ws.Range("I13").Value = "1"

$ws.Range("A14").Value = "6"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = "15 kΩ"
$ws.Range("D14").Value = "0402WGF1502TCE"
$ws.Range("E14").Value = "UNI-ROYAL(Uniroyal Elec)"
$ws.Range("F14").Value = "C25756"
$ws.Range("G14").Value = "R_0402_1005Metric"
$ws.Range("H14").Value = "62.5mW Thick Film Resistors ±1% 15kΩ 0402 Chip Resistor - Surface Mount ROHS"
$ws.Range("I14").Value = "6"

$ws.Range("A15").Value = "7"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "ATtiny10-TS"
$ws.Range("D15").Value = "ATTINY10-TSHR"
$ws.Range("E15").Value = "Microchip Tech"
$ws.Range("F15").Value = "C128438"
$ws.Range("G15").Value = "SOT-23-6"
$ws.Range("H15").Value = "AVR 32Byte 12MHz 4 SOT-23-6 Microcontroller Units (MCUs/MPUs/SOCs) ROHS"
$ws.Range("I15").Value = "1"

$ws.Range("A16").Value = "8"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "LM78M05_TO252"
$ws.Range("D16").Value = "L78M05ABDT-TR"
$ws.Range("E16").Value = "STMicroelectronics"
$ws.Range("F16").Value = "C58069"
$ws.Range("G16").Value = "TO-252-2"
$ws.Range("H16").Value = "2dB@(120Hz) 500mA Fixed 5V~5V Positive 35V TO-252-2(DPAK) Linear Voltage Regulators (LDO) ROHS"
$ws.Range("I16").Value = "1"

# ---------------------------------------------------------------------------
# 6. Reposition/resize the KiBot logo picture so it keeps anchoring
#    from A1 to col H (index 7) at the same pixel size, now that columns
#    are wider and some rows are taller.
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)

$widthPts = 0
for ($c = 1; $c -le 7; $c++) { $widthPts += $ws.Columns($c).Width }
$widthPts += 45.87582677165354

$heightPts = 0
for ($r = 1; $r -le 53; $r++) { $heightPts += $ws.Rows($r).Height }
$heightPts += 0.1258267716535433

$shp.Width = $widthPts
$shp.Height = $heightPts
